$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two new columns P and Q ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the formatting (bold font, border, centered alignment) from O1
# onto the two newly added header cells so they match the rest of the
# header row.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    # Swap the I/K values and the M/O values in each row.
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1

    # Fill in the two new columns P and Q with value 2.
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
